$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 0.31321845026989964
$ws.Range("B28").Value = 0.32081427931253648
$ws.Range("C28").Value = -0.0045967754223132616

$ws.Range("A29").Value = 0.31209067422019865
$ws.Range("B29").Value = 0.32844251853513751
$ws.Range("C29").Value = -0.0046058413287692127

$ws.Range("A30").Value = 0.32026659637766808
$ws.Range("B30").Value = 0.32026659637766808
$ws.Range("C30").Value = -0.0046088449302042582
